# Fruta / hortaliza, semanal
# Insert a new daily price-report row at row 615 (pushing the existing
# rows 615-715 down to 616-716), and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 615, shifting everything below
# it (rows 615-715) down by one -> new last row becomes 716.
$ws.Rows("615:615").Insert()

# Fill in the new row with the reported data.
$ws.Range("A615").Value = 8
$ws.Range("B615").Value = "Terminal La Palmera de La Serena"
$ws.Range("C615").Value = "Coquimbo"
$ws.Range("D615").Value = 45218
$ws.Range("E615").Value = 4
$ws.Range("F615").Value = 100114001
$ws.Range("G615").Value = "Papa"
$ws.Range("H615").Value = "Cardinal"
$ws.Range("I615").Value = "1a (cosecha)"
$ws.Range("J615").Value = 2000
$ws.Range("K615").Value = 25000
$ws.Range("L615").Value = 26000
$ws.Range("M615").Value = 25500
$ws.Range("N615").Value = "$/saco 25 kilos"
$ws.Range("O615").Value = "Provincia del Elquí"
$ws.Range("P615").Value = 1020
$ws.Range("Q615").Value = 25
$ws.Range("R615").Value = "Hortaliza"
